$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '49.667.56'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.62%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.637.68'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.02%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '113.22'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +2.41%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '324.10'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -1.33%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.529'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.78%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.546'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -2.50%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.91'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -2.01%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '19.86'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -4.43%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.01%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +1.18%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.33'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.55%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.049.86'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.09%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.630.88'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.07%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.861'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -2.41%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '49.562.58'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.72%  '
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -3.12%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -3.33%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -1.76%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0₃0948'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.61%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '270.29'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -3.90%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '69.00'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -5.56%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -2.30%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.35'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -1.18%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.02%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.37'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +4.35%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.23'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.139'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -3.26%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '35.13'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -4.05%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '49.66'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.26%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.89%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0816'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +2.33%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.21%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '19.05'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -3.62%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +3.99%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.76%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.23%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '126.92'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +2.59%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -1.36%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '22.40'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.53%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0325'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +3.20%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -3.22%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.061.10'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.21%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -3.51%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +6.67%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.15'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -6.70%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.95'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -1.40%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '59.18'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +1.51%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -3.18%  '
